# Apply scheduled profit-tracking updates across multiple sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2310.3076
$ws.Range("I19").Value = 1995.3334
$ws.Range("J19").Value = 2580.2856
$ws.Range("K19").Value = 1995.3334
$ws.Range("L19").Value = 2580.2856
$ws.Range("M19").Value = -1820.3334
$ws.Range("N19").Value = -2930.2856

$ws.Range("H41").Value = 255.04762
$ws.Range("I41").Value = 215.61539
$ws.Range("J41").Value = 319.125
$ws.Range("K41").Value = 215.61539
$ws.Range("L41").Value = 319.125
$ws.Range("M41").Value = 224.38461
$ws.Range("N41").Value = -1199.125

$ws.Range("H55").Value = 184.09091
$ws.Range("I55").Value = 194.2
$ws.Range("J55").Value = 175.66667
$ws.Range("K55").Value = 194.2
$ws.Range("L55").Value = 175.66667
$ws.Range("M55").Value = 19.80000000000001
$ws.Range("N55").Value = -603.6666700000001

$ws.Range("H97").Value = 477.5
$ws.Range("J97").Value = 477.5
$ws.Range("L97").Value = 1432.5
$ws.Range("N97").Value = -2424.5

$ws.Range("H112").Value = 1680.8485
$ws.Range("J112").Value = 1717.75
$ws.Range("L112").Value = 5153.25
$ws.Range("N112").Value = -7369.25

$ws.Range("H127").Value = 903.2222
$ws.Range("I127").Value = 697
$ws.Range("J127").Value = 911.1539
$ws.Range("K127").Value = 2091
$ws.Range("L127").Value = 2733.4617
$ws.Range("M127").Value = 2869
$ws.Range("N127").Value = -12653.4617

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").ClearContents()
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = 0

$ws.Range("H122").Value = 5613.517
$ws.Range("I122").Value = 6117.25
$ws.Range("K122").Value = 18351.75
$ws.Range("M122").Value = -15901.75

$ws.Range("H132").Value = 5103557
$ws.Range("I132").Value = 6098824
$ws.Range("J132").Value = 2816
$ws.Range("K132").Value = 18296472
$ws.Range("L132").Value = 8448
$ws.Range("M132").Value = -18293942
$ws.Range("N132").Value = -13508

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7756523.5
$ws.Range("I31").Value = 4751.933
$ws.Range("J31").Value = 25645226
$ws.Range("K31").Value = 4751.933
$ws.Range("L31").Value = 25645226
$ws.Range("M31").Value = -4456.933
$ws.Range("N31").Value = -25645816

$ws.Range("H34").Value = 7756523.5
$ws.Range("I34").Value = 4751.933
$ws.Range("J34").Value = 25645226
$ws.Range("K34").Value = 4751.933
$ws.Range("L34").Value = 25645226
$ws.Range("M34").Value = -4549.933
$ws.Range("N34").Value = -25645630

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").ClearContents()
$ws.Range("N63").Value = 0

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").ClearContents()
$ws.Range("N66").Value = 0

$ws.Range("H68").Value = 23899
$ws.Range("J68").Value = 23899
$ws.Range("L68").Value = 23899
$ws.Range("N68").Value = -25397

$ws.Range("H71").Value = 23899
$ws.Range("J71").Value = 23899
$ws.Range("L71").Value = 71697
$ws.Range("N71").Value = -79185

$ws.Range("H74").Value = 13714
$ws.Range("J74").Value = 13714
$ws.Range("L74").Value = 13714
$ws.Range("N74").Value = -15462

$ws.Range("H77").Value = 13714
$ws.Range("J77").Value = 13714
$ws.Range("L77").Value = 41142
$ws.Range("N77").Value = -49878

$ws.Range("H105").Value = 1104.8572
$ws.Range("I105").Value = 1162.5454
$ws.Range("K105").Value = 1162.5454
$ws.Range("M105").Value = 584.4546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 493.33334
$ws.Range("I39").Value = 30
$ws.Range("J39").Value = 520.58826
$ws.Range("K39").Value = 90
$ws.Range("L39").Value = 1561.76478
$ws.Range("M39").Value = 204
$ws.Range("N39").Value = -2149.76478

$ws.Range("H76").Value = 4295.8125
$ws.Range("I76").Value = 577.6667
$ws.Range("J76").Value = 5153.846
$ws.Range("K76").Value = 1733.0001
$ws.Range("L76").Value = 15461.538
$ws.Range("M76").Value = -1350.0001
$ws.Range("N76").Value = -16227.538

$ws.Range("H79").Value = 4295.8125
$ws.Range("I79").Value = 577.6667
$ws.Range("J79").Value = 5153.846
$ws.Range("K79").Value = 1733.0001
$ws.Range("L79").Value = 15461.538
$ws.Range("M79").Value = -407.0001
$ws.Range("N79").Value = -18113.538

$ws.Range("H82").Value = 6764.778
$ws.Range("I82").Value = 2506.5
$ws.Range("J82").Value = 7981.4287
$ws.Range("K82").Value = 7519.5
$ws.Range("L82").Value = 23944.2861
$ws.Range("M82").Value = -7113.5
$ws.Range("N82").Value = -24756.2861

$ws.Range("H85").Value = 6764.778
$ws.Range("I85").Value = 2506.5
$ws.Range("J85").Value = 7981.4287
$ws.Range("K85").Value = 7519.5
$ws.Range("L85").Value = 23944.2861
$ws.Range("M85").Value = -6115.5
$ws.Range("N85").Value = -26752.2861

$ws.Range("H88").Value = 3716.6667
$ws.Range("I88").Value = 3200
$ws.Range("J88").Value = 3975
$ws.Range("K88").Value = 9600
$ws.Range("L88").Value = 11925
$ws.Range("M88").Value = -9172
$ws.Range("N88").Value = -12781

$ws.Range("H91").Value = 3716.6667
$ws.Range("I91").Value = 3200
$ws.Range("J91").Value = 3975
$ws.Range("K91").Value = 9600
$ws.Range("L91").Value = 11925
$ws.Range("M91").Value = -8118
$ws.Range("N91").Value = -14889

$ws.Range("H118").Value = 1211
$ws.Range("J118").Value = 1385.4117
$ws.Range("L118").Value = 4156.2351
$ws.Range("N118").Value = -6642.2351

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10381.667
$ws.Range("I70").Value = 15658.706
$ws.Range("J70").Value = 4774.8125
$ws.Range("K70").Value = 15658.706
$ws.Range("L70").Value = 4774.8125
$ws.Range("M70").Value = -15388.706
$ws.Range("N70").Value = -5314.8125

$ws.Range("H73").Value = 10381.667
$ws.Range("I73").Value = 15658.706
$ws.Range("J73").Value = 4774.8125
$ws.Range("K73").Value = 15658.706
$ws.Range("L73").Value = 4774.8125
$ws.Range("M73").Value = -14722.706
$ws.Range("N73").Value = -6646.8125

$ws.Range("H126").Value = 3548.6155
$ws.Range("I126").Value = 2291.8333
$ws.Range("J126").Value = 4625.857
$ws.Range("K126").Value = 6875.499899999999
$ws.Range("L126").Value = 13877.571
$ws.Range("M126").Value = -4405.499899999999
$ws.Range("N126").Value = -18817.571

$ws.Range("H132").Value = 6626.4287
$ws.Range("I132").Value = 6651.6
$ws.Range("K132").Value = 19954.8
$ws.Range("M132").Value = -17424.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2263.7693
$ws.Range("I100").Value = 2300
$ws.Range("J100").Value = 2247.6667
$ws.Range("K100").Value = 2300
$ws.Range("L100").Value = 2247.6667
$ws.Range("M100").Value = -1759
$ws.Range("N100").Value = -3329.6667

$ws.Range("H132").Value = 7701681
$ws.Range("I132").Value = 7436.5
$ws.Range("J132").Value = 29426606
$ws.Range("K132").Value = 22309.5
$ws.Range("L132").Value = 88279818
$ws.Range("M132").Value = -19779.5
$ws.Range("N132").Value = -88284878

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5750.5
$ws.Range("I132").Value = 6502
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 19506
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -16976
$ws.Range("N132").Value = -20057
